# Change the deck's theme (Design tab) from "Integral" (Red Violet) to the
# standard "Office Theme" color scheme, matching the target colors found in
# ppt/theme/theme2.xml being promoted into ppt/theme/theme1.xml.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$tcs = $theme.ThemeColorScheme

# Office Theme colour scheme, in clrScheme order:
#  1 dk1      000000
#  2 lt1      FFFFFF
#  3 dk2      44546A
#  4 lt2      E7E6E6
#  5 accent1  5B9BD5
#  6 accent2  ED7D31
#  7 accent3  A5A5A5
#  8 accent4  FFC000
#  9 accent5  4472C4
# 10 accent6  70AD47
# 11 hlink    0563C1
# 12 folHlink 954F72
$tcs.Colors(1).RGB = 0
$tcs.Colors(2).RGB = 16777215
$tcs.Colors(3).RGB = 6968388
$tcs.Colors(4).RGB = 15132391
$tcs.Colors(5).RGB = 13998939
$tcs.Colors(6).RGB = 3243501
$tcs.Colors(7).RGB = 10855845
$tcs.Colors(8).RGB = 49407
$tcs.Colors(9).RGB = 12874308
$tcs.Colors(10).RGB = 4697456
$tcs.Colors(11).RGB = 12673797
$tcs.Colors(12).RGB = 7491477
